$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E) were changed to reflect a different
# sub-selection of subjects (16 / 20 / 16 / 20 instead of 1 / 2 / 3 / 4).
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) data for the same four columns.
$ws.Range("B2").Value = 11.851118303041785
$ws.Range("C2").Value = 11.371378553419259
$ws.Range("D2").Value = 12.782451313437109
$ws.Range("E2").Value = 11.910234114733127

# Row 3 (STR) data for the same four columns.
$ws.Range("B3").Value = 10.883366192557336
$ws.Range("C3").Value = 10.020953517009765
$ws.Range("D3").Value = 11.881610585812835
$ws.Range("E3").Value = 11.399252102828305

# Selection now only spans the four edited columns instead of the whole table.
[void]$ws.Range("B1:E3").Select()
